$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 1.65
$ws.Range("G5").Value = 1.42
$ws.Range("K5").Value = 5.6
$ws.Range("Q5").Value = 1.82
$ws.Range("T5").Value = 2.1
$ws.Range("X5").Value = 22
$ws.Range("AE5").Value = 150
$ws.Range("AK5").Value = 15.5
$ws.Range("AO5").Value = 220
$ws.Range("H6").Value = 1.8
$ws.Range("AF6").Value = 44
$ws.Range("AL6").Value = 85
$ws.Range("F7").Value = 1.22
$ws.Range("H7").Value = 15
$ws.Range("O7").Value = 1.13
$ws.Range("Q7").Value = 1.43
$ws.Range("R7").Value = 1.86
$ws.Range("S7").Value = 2.04
$ws.Range("T7").Value = 2
$ws.Range("X7").Value = 60
$ws.Range("AC7").Value = 22
$ws.Range("AK7").Value = 13.5
$ws.Range("F8").Value = 1.77
$ws.Range("H8").Value = 4.6
$ws.Range("N8").Value = 5.2
$ws.Range("O8").Value = 1.2
$ws.Range("F9").Value = 1.56
$ws.Range("G9").Value = 1.58
$ws.Range("P9").Value = 2.3
$ws.Range("Q9").Value = 1.68
$ws.Range("U9").Value = 2.14
$ws.Range("Z9").Value = 370
$ws.Range("AN9").Value = 7
$ws.Range("F10").Value = 1.51
$ws.Range("K10").Value = 5.1
$ws.Range("Q10").Value = 1.66
$ws.Range("Z10").Value = 760
$ws.Range("AC10").Value = 11
$ws.Range("H11").Value = 2.38
$ws.Range("I11").Value = 2.44
$ws.Range("M11").Value = 1.07
$ws.Range("X11").Value = 15
$ws.Range("G12").Value = 1.58
$ws.Range("H12").Value = 6.2
$ws.Range("X12").Value = 32
$ws.Range("AB12").Value = 10.5
$ws.Range("AK12").Value = 16
$ws.Range("I13").Value = 1.8
$ws.Range("H14").Value = 3.4
$ws.Range("J14").Value = 3.25
$ws.Range("Q14").Value = 1.86
$ws.Range("I17").Value = 3.15
$ws.Range("Q17").Value = 1.75
$ws.Range("G20").Value = 5.5
$ws.Range("H20").Value = 1.79
$ws.Range("I20").Value = 1.94
$ws.Range("J20").Value = 3.55
$ws.Range("P20").Value = 1.83
$ws.Range("Q20").Value = 1.99
$ws.Range("M22").Value = 1.04
$ws.Range("X22").Value = 32
$ws.Range("Q23").Value = 1.71
$ws.Range("AB23").Value = 10.5
$ws.Range("P24").Value = 2.42
$ws.Range("U24").Value = 2.52
$ws.Range("Z24").Value = 32
$ws.Range("AE24").Value = 44
$ws.Range("Q25").Value = 1.94
$ws.Range("G26").Value = 5.2
$ws.Range("H26").Value = 1.78
$ws.Range("I26").Value = 1.81
$ws.Range("P26").Value = 2
$ws.Range("Q26").Value = 1.88
$ws.Range("T26").Value = 1.82
$ws.Range("X26").Value = 17
$ws.Range("Z26").Value = 11.5
$ws.Range("AI26").Value = 1000
$ws.Range("AN26").Value = 1000
$ws.Range("G27").Value = 1.4
$ws.Range("I27").Value = 10.5
$ws.Range("J27").Value = 5.4
$ws.Range("K27").Value = 5.9
$ws.Range("M27").Value = 1.04
$ws.Range("P27").Value = 2.3
$ws.Range("T27").Value = 2.04
$ws.Range("U27").Value = 1.87
$ws.Range("X27").Value = 32
$ws.Range("Y27").Value = 34
$ws.Range("Z27").Value = 1000
$ws.Range("AA27").Value = 440
$ws.Range("AD27").Value = 65
$ws.Range("AF27").Value = 8.4
$ws.Range("AH27").Value = 50
$ws.Range("AI27").Value = 1000
$ws.Range("AJ27").Value = 11.5
$ws.Range("AK27").Value = 15.5
$ws.Range("AL27").Value = 40
$ws.Range("AN27").Value = 5.9
$ws.Range("H28").Value = 7.4
$ws.Range("Q28").Value = 1.75
$ws.Range("T28").Value = 1.9
$ws.Range("AL28").Value = 75
$ws.Range("H29").Value = 15
$ws.Range("M29").Value = 1.03
$ws.Range("N29").Value = 5.4
$ws.Range("O29").Value = 1.2
$ws.Range("P29").Value = 2.54
$ws.Range("Q29").Value = 1.58
$ws.Range("R29").Value = 1.65
$ws.Range("T29").Value = 2.28
$ws.Range("X29").Value = 28
$ws.Range("AB29").Value = 9.4
$ws.Range("AC29").Value = 17.5
$ws.Range("AD29").Value = 65
$ws.Range("AE29").Value = 350
$ws.Range("AG29").Value = 12.5
$ws.Range("J31").Value = 3.45
$ws.Range("P31").Value = 2.02
$ws.Range("F33").Value = 1.45
$ws.Range("J33").Value = 3.8
$ws.Range("K33").Value = 5.2
$ws.Range("J34").Value = 3.6
$ws.Range("F36").Value = 6.2
$ws.Range("G36").Value = 8.4
$ws.Range("H36").Value = 1.43
$ws.Range("F37").Value = 1.71
$ws.Range("I37").Value = 6.2
$ws.Range("Q37").Value = 2.04
$ws.Range("AB38").Value = 15.5
$ws.Range("AD38").Value = 13
$ws.Range("AF38").Value = 28
$ws.Range("AG38").Value = 17
$ws.Range("AH38").Value = 20
$ws.Range("AJ38").Value = 65
$ws.Range("AL38").Value = 50
$ws.Range("AN38").Value = 44
$ws.Range("AO38").Value = 22
$ws.Range("F40").Value = 1.86
